# Apply changes to the "247 stress-perm" sheet:
#  - Row 16: clear "Running?" (G16), update Remarks (I16) text, add Results (J16) text
#  - Add new rows 17-21 (Subcases 14-18) cloning row 16's pattern with updated Remarks/Results
#  - Leave row 22 blank
#  - Add new rows 23-25 (Subcases 20-22) with a different Remarks/Results pair, plus
#    "Saved as" on row 25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("247 stress-perm")

$xlHAlignCenter = -4108

# --- Update existing row 16 ---
$ws.Cells.Item(16, 7).Value = $null
$ws.Cells.Item(16, 10).Value = "better, but validation is problem"
$ws.Cells.Item(16, 9).Value = "modified again CNN10, and more changes"

# Reference formatting pulled from row 16, which already carries the look
# the new rows should copy (centered numbers/booleans in D:F, red-font
# centered "Running?" cell in G).
$dAlign = $ws.Cells.Item(16, 4).HorizontalAlignment
$gColor = $ws.Cells.Item(16, 7).Font.Color

function Format-NewRow($row) {
    $ws.Range($ws.Cells.Item($row, 4), $ws.Cells.Item($row, 6)).HorizontalAlignment = $dAlign
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.HorizontalAlignment = $dAlign
    $gCell.Font.Color = $gColor
}

# --- New rows 17-21 (same Job id/Name/CNN, Subcases 14..18) ---
$subcases = 14,15,16,17,18
$row = 17
foreach ($sc in $subcases) {
    $ws.Cells.Item($row, 2).Value = "my laptop"
    $ws.Cells.Item($row, 3).Value = "recognizeStressPermf"
    $ws.Cells.Item($row, 4).Value = 10
    $ws.Cells.Item($row, 5).Value = $sc
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 9).Value = "modified again CNN10, and more changes"
    $ws.Cells.Item($row, 10).Value = "better, but validation is problem"
    Format-NewRow $row
    $row++
}

# Row 22 intentionally left blank

# --- New rows 23-25 (Subcases 20..22) ---
$subcases2 = 20,21,22
$row = 23
foreach ($sc in $subcases2) {
    $ws.Cells.Item($row, 2).Value = "my laptop"
    $ws.Cells.Item($row, 3).Value = "recognizeStressPermf"
    $ws.Cells.Item($row, 4).Value = 10
    $ws.Cells.Item($row, 5).Value = $sc
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 9).Value = "also  used manual train/valid splitting"
    $ws.Cells.Item($row, 10).Value = "same as above"
    Format-NewRow $row
    $row++
}

# "Saved as" for the last new row
$ws.Cells.Item(25, 11).Value = "saved"

# Match column widths to the new content (Excel auto-fit after the longer text
# was entered in columns H, I, J, and a slight narrowing of G). The inputs
# below are chosen so the engine's internally-quantized stored width lands
# on the value closest to the target column widths.
$ws.Columns.Item(7).ColumnWidth = 9.165
$ws.Columns.Item(8).ColumnWidth = 11.332
$ws.Columns.Item(9).ColumnWidth = 35.003
$ws.Columns.Item(10).ColumnWidth = 28.332

# Restore the selection to where the author left off.
$ws.Range("J28").Select()
